$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Betarraga" (Macroferia Regional de Talca) needs
# to be inserted at row 134, pushing the existing rows 134-174 down to 135-175.
$ws.Rows(134).Insert()

# Fill in the newly inserted row 134 with the new record's data.
$ws.Cells.Item(134, 1).Value = 5
$ws.Cells.Item(134, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(134, 3).Value = "Maule"
$ws.Cells.Item(134, 4).Value = 44463
$ws.Cells.Item(134, 5).Value = 7
$ws.Cells.Item(134, 6).Value = 100114014
$ws.Cells.Item(134, 7).Value = "Betarraga"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 3000
$ws.Cells.Item(134, 11).Value = 650
$ws.Cells.Item(134, 12).Value = 650
$ws.Cells.Item(134, 13).Value = 650
$ws.Cells.Item(134, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(134, 15).Value = "Región del Maule"
$ws.Cells.Item(134, 16).Value = 130
$ws.Cells.Item(134, 17).Value = 5
$ws.Cells.Item(134, 18).Value = "Hortaliza"
